$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Two new data rows (26 / 27) for the new "主世界bgm1" / "主世界bgm2"
#    sound entries, following the same pattern as the existing rows.
# ------------------------------------------------------------------

# Column A (id) + D..I keep the same look-and-feel as the row above
# (row 29) -- copy its number/left-aligned formatting down first.
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)

$ws.Range("D29:I29").Copy()
$ws.Range("D30:I30").PasteSpecial(-4122)
$ws.Range("D31:I31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A30").Value = 26
$ws.Range("C30").Value = "主世界bgm1"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0.4
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 200
$ws.Range("I30").Value = 600

$ws.Range("A31").Value = 27
$ws.Range("C31").Value = "主世界bgm2"
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0.4
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 200
$ws.Range("I31").Value = 600

# ------------------------------------------------------------------
# 2) Trailing blank rows 32..47 -- only column A carries the row
#    formatting (mirrors dragging the fill handle down with no data).
# ------------------------------------------------------------------
$ws.Range("A29").Copy()
$ws.Range("A32:A47").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) View state: scroll + selection as left by the author.
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J8").Select()
